$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to be stored as text (preserving leading zeros, e.g.
    # account numbers like "005105172"), then strip the temporary
    # NumberFormat override so the cell ends up unstyled again, matching
    # the other plain inlineStr cells in this sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# 1) Insert a new row for account 005105172 / VALDIVINO / 70101.54
#    just above the existing "004752534 / CARLOS" row (currently row 3).
$ws.Rows.Item(3).Insert()
Set-TextCell $ws.Cells.Item(3, 1) "005105172"
$ws.Cells.Item(3, 2).Value = "VALDIVINO"
$ws.Cells.Item(3, 3).Value = 70101.54

# 2) Insert a new row for account 004458563 / LUIZ / 61003.79
#    just above the existing "005101676 / ELENI" row (now row 5 after step 1).
$ws.Rows.Item(5).Insert()
Set-TextCell $ws.Cells.Item(5, 1) "004458563"
$ws.Cells.Item(5, 2).Value = "LUIZ"
$ws.Cells.Item(5, 3).Value = 61003.79

# 3) The row that used to hold 005105172 / VALDIVINO / 42029.42 has shifted
#    down by two rows (from 7 to 9); overwrite it with 004575632 / ADELE / 41063.96.
Set-TextCell $ws.Cells.Item(9, 1) "004575632"
$ws.Cells.Item(9, 2).Value = "ADELE"
$ws.Cells.Item(9, 3).Value = 41063.96

# 4) Remove the now-duplicate row that held 004575632 / ADELE / 23566.46
#    (originally row 9, shifted down by two rows to row 11).
$ws.Rows.Item(11).Delete()
